# Habitat_Quality_and_Geomorphic_Potential_Rating_Criteria.xlsx edit
# Updates the Geomorphic Potential "unconfined" percentage thresholds to
# finalized percentage-based categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: <30% unconfined  ->  <0.5%   (upper bound 30 -> 0.5)
$ws.Range("H26").Value = 0.5
$ws.Range("J26").Value = "<0.5%"

# Row 27: 30-60% unconfined -> 0.5%-60%  (lower bound 30 -> 0.5)
$ws.Range("G27").Value = 0.5
$ws.Range("J27").Value = "0.5%-60%"

# Freeze pane / scroll position and selection moved further down the sheet
$ws.Activate()
$window = $excel.ActiveWindow
$window.SplitRow = 1
$window.FreezePanes = $true
$ws.Range("A14").Select()
$window.ScrollRow = 14

$ws.Range("D1").Select()
$ws.Range("C29").Select()

# Window size/position on save
$window.WindowState = -4143
$window.Top = 1536
$window.Left = 1536
$window.Width = 21156
$window.Height = 11412
